$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) — rows 3,4,8,12,16,18,19,21,22
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 391
$ws1.Range("F4").Value  = 1496
$ws1.Range("F8").Value  = 646
$ws1.Range("F12").Value = 3589
$ws1.Range("F16").Value = 1184
$ws1.Range("F18").Value = 1122
$ws1.Range("F19").Value = 308
$ws1.Range("F21").Value = 2340
$ws1.Range("F22").Value = 56

# Sheet "全部类型" (all types) — same events, but the "银魂主题派对" row is shifted
# down by one (row 23 instead of row 22) because this sheet also includes the
# "演出" (performance) event inserted at row 22.
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 391
$ws4.Range("F4").Value  = 1496
$ws4.Range("F8").Value  = 646
$ws4.Range("F12").Value = 3589
$ws4.Range("F16").Value = 1184
$ws4.Range("F18").Value = 1122
$ws4.Range("F19").Value = 308
$ws4.Range("F21").Value = 2340
$ws4.Range("F23").Value = 56
